$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray _GoBack bookmark that currently sits in the title
#    paragraph (it will be re-created later, anchored to "DX Grid").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Locate the two list-item paragraphs that describe double-clicking a
#    decoded message / callsign. They get merged into a single bullet.
# ---------------------------------------------------------------------------
$firstIdx = -1
$secondIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($firstIdx -eq -1 -and $t -match "^Click on waterfall to set QSO Frequency\.") {
        $firstIdx = $i
    }
    if ($secondIdx -eq -1 -and $t -match "^Double-click on a decoded callsign") {
        $secondIdx = $i
    }
}

if ($firstIdx -eq -1 -or $secondIdx -eq -1) {
    throw "Could not locate target paragraphs (first=$firstIdx, second=$secondIdx)"
}

# ---------------------------------------------------------------------------
# 3. Rewrite the first paragraph's text completely (this also drops the
#    second paragraph's distinctive lead-in, since its content is folded
#    into this single bullet).
# ---------------------------------------------------------------------------
$r1 = $d.Paragraphs($firstIdx).Range
$r1.MoveEnd(1, -1) | Out-Null
$r1.Text = "Double-click on a decoded message to copy the callsign and locator into DX Call and DX Grid.  This will also generate appropriate standard messages, including signal report."

# ---------------------------------------------------------------------------
# 4. Delete the now-redundant second paragraph entirely (its wording has
#    been absorbed into the rewritten first paragraph above).
# ---------------------------------------------------------------------------
$d.Paragraphs($secondIdx).Range.Delete()

# ---------------------------------------------------------------------------
# 5. Re-apply bold character formatting to "DX Call" and "DX Grid", and wrap
#    "DX Grid" with a fresh _GoBack bookmark, matching the original markup
#    style used for such field-name callouts.
# ---------------------------------------------------------------------------
$rerange = $d.Paragraphs($firstIdx).Range
$foundCall = $rerange.Find.Execute("DX Call", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundCall) {
    $rerange.Bold = 1
}

$rerange2 = $d.Paragraphs($firstIdx).Range
$foundGrid = $rerange2.Find.Execute("DX Grid", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundGrid) {
    $rerange2.Bold = 1
    $d.Bookmarks.Add("_GoBack", $rerange2)
}

# ---------------------------------------------------------------------------
# 6. The footer's cached PAGE field result needs to reflect the document's
#    new (shorter) page count: 6 -> 5.
# ---------------------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$ffind = $footer.Range
$ffind.Find.Execute("6", $false, $false, $false, $false, $false, $true, 1, $false, "5", 2) | Out-Null
